# Variable_Overview.xlsx — "Add external variables to pilot"
#
# Adds two new variable-table rows to the Pilot sheet (type_initialization,
# time_initialization) plus a blank spacer row, right after the header row,
# pushing the existing variable rows down. Also leaves the Pilot sheet as
# the active/selected sheet (matching the saved view state in the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pilot")

# --- Stash the formatting of the two rows that are about to shift down so
#     the freshly inserted rows can reuse the exact same cell styles
#     (s="6" header-like style currently on row 5, s="5" body style on row 6).
#     We stash into a scratch area far to the right (columns H:L) and wipe
#     it again once the real rows have their formats back.
$ws.Range("B5:F5").Copy()
$ws.Range("H1:L1").PasteSpecial(-4122)
$ws.Range("B6:F6").Copy()
$ws.Range("H2:L2").PasteSpecial(-4122)

# --- Insert 3 new rows above the old row 5 (2 data rows + 1 blank spacer),
#     shifting everything below down (old row 5 -> row 8, ..., old row 27 -> row 30).
$ws.Rows("5:7").Insert()

# --- Restore styles onto the newly inserted rows.
$ws.Range("H1:L1").Copy()
$ws.Range("B5:F5").PasteSpecial(-4122)
$ws.Range("H2:L2").Copy()
$ws.Range("B6:F6").PasteSpecial(-4122)
$ws.Range("B7:F7").PasteSpecial(-4122)

# --- Clean up the scratch area.
$ws.Range("H1:L2").Clear()

# --- Fill in the new row content. Order matches how the values were
#     originally authored (row 5's Type/Unit/Value Range, then all of row 6,
#     then row 5's Description/Variable Name last) so shared strings line up.
$ws.Range("C5").Value = "string"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = "[Landing, Takeoff]"

$ws.Range("B6").Value = "time_initialization"
$ws.Range("C6").Value = "integer"
$ws.Range("D6").Value = "s"
$ws.Range("E6").Value = "[2; inf]"
$ws.Range("F6").Value = "Defines"

$ws.Range("F5").Value = "Defines whether this pilot is going to be landing or starting from ground"
$ws.Range("B5").Value = "type_initialization"

# Row 5 wraps onto two lines (matches the new header-like row's taller height).
$ws.Rows(5).RowHeight = 28.5

# --- Leave the Pilot sheet active/selected (it becomes the saved active tab).
$ws.Activate()
$ws.Range("B5").Select()
